# Template_Capex_export.xlsx — "Fix bug capex, export by dept"
#
# The CAPEX sheet had two columns pre-seeded with placeholder defaults for
# every data row (QTY = 1 in column N, Price = 0 in column R). That caused
# the per-department export/budget calculations to pick up bogus 1/0
# values for rows the user never actually filled in. The fix clears those
# placeholder values so the cells are genuinely blank until the user
# enters real data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAPEX")

# Data rows run from row 4 through row 52 (the "No" column counts 3..51).
# Column N = QTY (currently pre-filled with 1), column R = Price
# (currently pre-filled with 0). Clear both so the cells are blank.
for ($r = 4; $r -le 52; $r++) {
    $ws.Cells.Item($r, 14).ClearContents()   # N<r> — QTY
    $ws.Cells.Item($r, 18).ClearContents()   # R<r> — Price
}

# Refresh the saved view state: scroll the grid so column M is the
# left-most visible column, and leave the selection on S6.
$win = $excel.ActiveWindow
$win.ScrollColumn = 13
$win.ScrollRow = 1
$ws.Range("S6").Select()

$wb.Save()
